$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registro")
$ws.Range("A9:P9").Insert()
Write-Host "inserted"
